# João Rodrigues schedule fix: shift afternoon block by one slot (6 hours by turn fix)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (7:00) ---
$ws.Range("F2").Value = "-"

# --- Row 3 (7:50) ---
$ws.Range("D3").Value = "MCT-2A-CAD"
$ws.Range("E3").Value = "MEC-1A-Desenho Técnico"

# --- Row 4 (8:40) ---
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"

# --- Row 6 (9:50) ---
$ws.Range("D6").Value = "ELT-1A-Desenho Técnico"
$ws.Range("E6").Value = "MEC-2A-CAD"
$ws.Range("F6").Value = "MEC-2A-CAD"

# --- Row 7 (10:40) ---
$ws.Range("D7").Value = "ELT-1A-Desenho Técnico"
$ws.Range("E7").Value = "MCT-2A-CAD"

# --- Row 8 (11:30) : Almoço moves out of this row ---
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "-"

# --- Row 9 : time becomes 12:20, Almoço moves in ---
$ws.Range("A9").Value = "12:20"
$ws.Range("B9").Value = "Almoço"
$ws.Range("C9").Value = "Almoço"
$ws.Range("D9").Value = "Almoço"
$ws.Range("E9").Value = "Almoço"
$ws.Range("F9").Value = "Almoço"

# --- Row 10 : time shifts from 13:50 to 13:00 ---
$ws.Range("A10").Value = "13:00"

# --- Row 11 : time shifts from 14:40 to 13:50 ---
$ws.Range("A11").Value = "13:50"

# --- Row 12 : time shifts from 15:30 to 14:40, Intervalo moves out ---
$ws.Range("A12").Value = "14:40"
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "-"

# --- Row 13 : time shifts from 15:50 to 15:30, Intervalo moves in ---
$ws.Range("A13").Value = "15:30"
$ws.Range("B13").Value = "Intervalo"
$ws.Range("C13").Value = "Intervalo"
$ws.Range("D13").Value = "Intervalo"
$ws.Range("E13").Value = "Intervalo"
$ws.Range("F13").Value = "Intervalo"

# --- Row 14 (new) : 15:50, all '-' ---
$ws.Range("A14").Value = "15:50"
$ws.Range("B14").Value = "-"
$ws.Range("C14").Value = "-"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "-"

# --- Row 15 (new) : 16:40, all '-' ---
$ws.Range("A15").Value = "16:40"
$ws.Range("B15").Value = "-"
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"

# --- Row 16 (new) : 17:30, all '-' ---
$ws.Range("A16").Value = "17:30"
$ws.Range("B16").Value = "-"
$ws.Range("C16").Value = "-"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "-"

# --- Row 17 (new) : 18:20, rest empty (but still present as blank cells) ---
$ws.Range("A17").Value = "18:20"
$ws.Range("B17:F17").Borders.LineStyle = 0

$wb.Save()
